{"js": "// Add bill Payment backend:\n// Remove the stray/duplicated \"DateEdit\" bullet describing the soft-delete\n// timestamp for \"danh m\u1ee5c\" (category) \u2014 this sentence was mistakenly\n// duplicated under the Category section and is removed here.\nconst body = context.document.body;\n\nconst searchText =\n  \"DateEdit: Th\u1ec3 hi\u1ec7n th\u1eddi gian x\u00f3a m\u1ec1m c\u1ee7a danh m\u1ee5c v\u1edbi m\u1ee5c \u0111\u00ednh d\u00e0nh cho vi\u1ec7c t\u1ef1 \u0111\u1ed9ng x\u00f3a c\u1ee9ng sau m\u1ed9t kho\u1ea3ng th\u1eddi gian.\";\n\nconst results = body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Delete the whole paragraph that contains the matched sentence so no\n  // empty bullet line is left behind.\n  const paragraph = results.items[0].paragraphs.getFirst();\n  paragraph.delete();\n  await context.sync();\n}\n", "ps1": "# Add bill Payment backend:\n# Remove the stray/duplicated \"DateEdit\" bullet describing the soft-delete\n# timestamp for \"danh m\u1ee5c\" (category) \u2014 this sentence was mistakenly\n# duplicated under the Category section and is removed here.\n\n$d = $word.ActiveDocument\n\n$searchText = \"DateEdit: Th\u1ec3 hi\u1ec7n th\u1eddi gian x\u00f3a m\u1ec1m c\u1ee7a danh m\u1ee5c v\u1edbi m\u1ee5c \u0111\u00ednh d\u00e0nh cho vi\u1ec7c t\u1ef1 \u0111\u1ed9ng x\u00f3a c\u1ee9ng sau m\u1ed9t kho\u1ea3ng th\u1eddi gian.\"\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.MatchCase = $true\n$found = $range.Find.Execute($searchText)\n\nif ($found) {\n    # $range now covers the matched sentence; expand to the enclosing\n    # paragraph and delete it (including its paragraph mark) so no blank\n    # bullet line is left behind.\n    $para = $range.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
